# Update news.xlsx rows 2-13: Title (C), Link (D) and Extraction Date (E)
# The "Site" column (B) and the index column (A) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Title = "PF faz operação contra a Braskem após afundamento de solo em Maceió"; Link = "https://g1.globo.com/al/alagoas/noticia/2023/12/21/pf-cumpre-mandados-em-investigacao-que-apura-crimes-cometidos-pela-braskem-em-maceio.ghtml"; Date = "21/12/2023 09:40:13" }
    @{ Row = 3;  Title = "Brasil reverteu 'veto' de Israel a 16 brasileiros e parentes"; Link = "https://g1.globo.com/politica/blog/camila-bomfim/post/2023/12/21/metade-do-grupo-de-brasileiros-a-ser-repatriado-de-gaza-tinha-recebido-veto-de-israel-brasil-reverteu.ghtml"; Date = "21/12/2023 09:40:13" }
    @{ Row = 4;  Title = "Milei anuncia decreto com 300 medidas para desregular economia"; Link = "https://g1.globo.com/mundo/noticia/2023/12/20/milei-anuncia-decreto-que-estabelece-bases-para-novo-plano-economico.ghtml"; Date = "21/12/2023 09:40:13" }
    @{ Row = 5;  Title = "Exploração de sal-gema pela Braskem em Maceió é alvo de operação da PF"; Link = "https://noticias.uol.com.br/cotidiano/ultimas-noticias/2023/12/21/operacao-pf-maceio.htm"; Date = "21/12/2023 09:40:13" }
    @{ Row = 6;  Title = "Deputado chora após levar tapa na cara de colega do PT: 'Humilhado'"; Link = "https://noticias.uol.com.br/politica/ultimas-noticias/2023/12/20/messias-donato-pronunciamento.htm"; Date = "21/12/2023 09:40:13" }
    @{ Row = 7;  Title = "Após Javier Milei anunciar megadecreto, argentinos protestam em Buenos Aires"; Link = "https://noticias.uol.com.br/internacional/ultimas-noticias/2023/12/21/apos-anuncio-de-megadecreto-argentinos-protestam-em-buenos-aires.htm"; Date = "21/12/2023 09:40:13" }
    @{ Row = 8;  Title = "Pescador captura espécie de peixe invasora que ameaça biodiversidade"; Link = "https://www.terra.com.br/planeta/noticias/pescador-captura-especie-de-peixe-invasora-que-ameaca-biodiversidade-no-rn,cb581c9e42757b17caf9bfc2d5c875346j2516sb.html"; Date = "21/12/2023 09:40:13" }
    @{ Row = 9;  Title = "Homem fatura R$ 25 milhões ajudando pessoas a comprar o presente certo"; Link = "https://forbes.com.br/forbes-money/2023/12/ele-faturou-r-25-milhoes-em-2023-ajudando-pessoas-a-comprar-o-presente-certo/?utm_source=terra_capa_noticias&utm_medium=referral"; Date = "21/12/2023 09:40:13" }
    @{ Row = 10; Title = "Governo lança nova fase de operação de repatriação de brasileiros em Gaza"; Link = "https://www.terra.com.br/noticias/governo-lanca-nova-fase-de-operacao-de-repatriacao-de-brasileiros-em-gaza,5b217a63bc3016beb3c2d2e56f8e20b9xuwdzc7x.html"; Date = "21/12/2023 09:40:13" }
    @{ Row = 11; Title = "PF apreende quase 3 toneladas de drogas no Aeroporto de Guarulhos em 2023"; Link = "https://www.cnnbrasil.com.br/nacional/pf-apreende-quase-3-toneladas-de-drogas-no-aeroporto-de-guarulhos-em-2023/"; Date = "21/12/2023 09:40:13" }
    @{ Row = 12; Title = "Presidente do Corinthians sobre Gabigol: “Chegaremos a um acordo”"; Link = "https://www.cnnbrasil.com.br/esportes/presidente-do-corinthians-sobre-gabigol-chegaremos-a-um-acordo/"; Date = "21/12/2023 09:40:13" }
    @{ Row = 13; Title = "Programas de desenvolvimento humano elevam a eficiência empresarial"; Link = "https://www.cnnbrasil.com.br/branded-content/nacional/programas-de-desenvolvimento-humano-elevam-a-eficiencia-empresarial/"; Date = "21/12/2023 09:40:13" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Title
    $ws.Cells.Item($r, 4).Value = $item.Link
    $ws.Cells.Item($r, 5).Value = $item.Date
}
